$d = $word.ActiveDocument

function FindParagraphIndex($substr) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $p = $d.Paragraphs.Item($i)
        if ($p.Range.Text.Contains($substr)) {
            return $i
        }
    }
    return -1
}

function ReplaceSelf($findText) {
    $r = $d.Content
    $ok = $r.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, $findText, 2)
    return $ok
}

# -----------------------------------------------------------------
# 1. Move "Use cases needed to be completed" paragraph from its
#    original spot (right before "Search through rad document...")
#    down to just after the "...and delete" (watchlist) paragraph,
#    right before "Edit user settings Use Case:".
# -----------------------------------------------------------------
$srcIdx = FindParagraphIndex("Use cases needed to be completed")
$src = $d.Paragraphs.Item($srcIdx)
$src.Range.Cut()

$dstIdx = FindParagraphIndex("user observes my watchlist")
$dst = $d.Paragraphs.Item($dstIdx)
$endPos = $dst.Range.End
$insertPoint = $d.Range($endPos, $endPos)
$insertPoint.Paste()

# -----------------------------------------------------------------
# 2. Run-merges (adjacent same-text runs collapsing into one run) --
#    these happen naturally in Word when you re-run Find/Replace
#    with identical find/replace text over a region spanning
#    multiple same-meaning runs.
# -----------------------------------------------------------------
ReplaceSelf("Exception edit (") | Out-Null
ReplaceSelf("AD - Start Application Use Case:") | Out-Null
ReplaceSelf("AD - Logout Use Case:") | Out-Null
ReplaceSelf("AD - Reset Password Use Case: ") | Out-Null
ReplaceSelf("AD - View Portfolio Summary Use Case") | Out-Null
ReplaceSelf("Scenario edit 1") | Out-Null
ReplaceSelf(": " + [char]0x201C + "if there are not four items in the portfolio display message to the user (number) " + [char]0x201C + "add items to your portfolio" + [char]0x201D + " ") | Out-Null

# -----------------------------------------------------------------
# 3. Re-create the _GoBack bookmark at its new location: its own
#    empty paragraph right before the relocated
#    "Use cases needed to be completed" paragraph.
# -----------------------------------------------------------------
$dstIdx2 = FindParagraphIndex("user observes my watchlist")
$dst2 = $d.Paragraphs.Item($dstIdx2)
$dst2.Range.InsertParagraphAfter()
$newP = $d.Paragraphs.Item($dstIdx2 + 1)
$xmlFrag = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>'
$newP.Range.InsertXML($xmlFrag)

Write-Host "done"
